$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# Update values for Luty (row 3), Sierpień (row 9), Wrzesień (row 10), Listopad (row 12)
$ws.Range("B3").Value = 160
$ws.Range("C3").Value = 20

$ws.Range("B9").Value = 160
$ws.Range("C9").Value = 20
$ws.Range("D9").Value = 11

$ws.Range("B10").Value = 176
$ws.Range("C10").Value = 22
$ws.Range("D10").Value = 8

$ws.Range("B12").Value = 144
$ws.Range("C12").Value = 18
$ws.Range("D12").Value = 12

# Set page setup to portrait orientation
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

# Set selection to G18
$ws.Range("G18").Select() | Out-Null
